# Generate Report for Handback
#
# Refreshes the handoff/handback timestamps for the
# 56f5903a-27d0-450b-b254-6d866f433341 file (row 3 on each status sheet)
# after a new xliff round-trip, mirroring what the handback-status report
# generator writes out.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# "Latest HO Xliff Generate Date" for 56f5903a-...md
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-31 00:51:46"

# --- zh-cn sheet ------------------------------------------------------
# "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-31 00:51:41"
$zhcn.Range("K3").Value = "2016-08-31 00:51:58"

# --- de-de sheet ------------------------------------------------------
# "Correspond Handoff Datetime" / "Correspond Handback DateTime"
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-31 00:51:46"
$dede.Range("K3").Value = "2016-08-31 00:52:12"
